# Daily attendance processing - 2025-12-08 17:53:40
# Fix the ordering of names in the "Recorded By" column (G) so that the
# automated "System" entry is listed last instead of first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
}
